# Scheduled runner update: refresh Garuda Profits market-price snapshot cells.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 478.51428
$ws.Range("J17").Value = 529.6774
$ws.Range("L17").Value = 1589.0322
$ws.Range("N17").Value = -1925.0322

$ws.Range("H43").Value = 233335630
$ws.Range("I43").Value = 333335170
$ws.Range("J43").Value = 83336340
$ws.Range("K43").Value = 333335170
$ws.Range("L43").Value = 83336340
$ws.Range("M43").Value = -333335101
$ws.Range("N43").Value = -83336478

$ws.Range("H107").Value = 585.05554
$ws.Range("I107").Value = 459.23077
$ws.Range("J107").Value = 912.2
$ws.Range("K107").Value = 459.23077
$ws.Range("L107").Value = 912.2
$ws.Range("M107").Value = 1460.76923
$ws.Range("N107").Value = -4752.2

$ws.Range("H125").Value = 2228.5715
$ws.Range("I125").Value = 1900
$ws.Range("J125").Value = 2360
$ws.Range("K125").Value = 17100
$ws.Range("L125").Value = 21240
$ws.Range("M125").Value = -14640
$ws.Range("N125").Value = -26160

$ws.Range("H137").Value = 1446.1923
$ws.Range("I137").Value = 1229.0952
$ws.Range("J137").Value = 2358
$ws.Range("K137").Value = 3687.2856
$ws.Range("L137").Value = 7074
$ws.Range("M137").Value = -1137.2856
$ws.Range("N137").Value = -12174

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1073.6
$ws.Range("I45").Value = 1017
$ws.Range("J45").Value = 1300
$ws.Range("K45").Value = 1017
$ws.Range("L45").Value = 1300
$ws.Range("M45").Value = -640
$ws.Range("N45").Value = -2054

$ws.Range("H61").Value = 2125.647
$ws.Range("I61").Value = 1117
$ws.Range("J61").Value = 3022.2222
$ws.Range("K61").Value = 1117
$ws.Range("L61").Value = 3022.2222
$ws.Range("M61").Value = -905
$ws.Range("N61").Value = -3446.2222

$ws.Range("H63").Value = 771908.0600000001
$ws.Range("I63").Value = 836067.0600000001
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 836067.0600000001
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -835381.0600000001
$ws.Range("N63").Value = -3372

$ws.Range("H66").Value = 771908.0600000001
$ws.Range("I66").Value = 836067.0600000001
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 4180335.3
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -4176903.3
$ws.Range("N66").Value = -16864

$ws.Range("H88").Value = 560612.2
$ws.Range("I88").Value = 1004630.5
$ws.Range("J88").Value = 5589.25
$ws.Range("K88").Value = 1004630.5
$ws.Range("L88").Value = 5589.25
$ws.Range("M88").Value = -1004224.5
$ws.Range("N88").Value = -6401.25

$ws.Range("H91").Value = 560612.2
$ws.Range("I91").Value = 1004630.5
$ws.Range("J91").Value = 5589.25
$ws.Range("K91").Value = 1004630.5
$ws.Range("L91").Value = 5589.25
$ws.Range("M91").Value = -1003226.5
$ws.Range("N91").Value = -8397.25

$ws.Range("H132").Value = 9502.286
$ws.Range("I132").Value = 12003.777
$ws.Range("K132").Value = 36011.331
$ws.Range("M132").Value = -33481.331

$ws.Range("H136").Value = 2125.647
$ws.Range("I136").Value = 1117
$ws.Range("J136").Value = 3022.2222
$ws.Range("K136").Value = 3351
$ws.Range("L136").Value = 9066.6666
$ws.Range("M136").Value = -801
$ws.Range("N136").Value = -14166.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 392.30768
$ws.Range("I22").Value = 392.30768
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 392.30768
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -219.30768
$ws.Range("N22").ClearContents()

$ws.Range("H24").Value = 6225
$ws.Range("I24").Value = 950
$ws.Range("J24").Value = 11500
$ws.Range("K24").Value = 950
$ws.Range("L24").Value = 11500
$ws.Range("M24").Value = -715
$ws.Range("N24").Value = -11970

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 66670170
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 83336960
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 83336960
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -83338208

$ws.Range("H65").Value = 66670170
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 83336960
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 416684800
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -416691040

$ws.Range("H134").Value = 1165.75
$ws.Range("I134").Value = 1165.75
$ws.Range("K134").Value = 3497.25
$ws.Range("M134").Value = -962.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 915.375
$ws.Range("I129").Value = 798
$ws.Range("J129").Value = 1111
$ws.Range("K129").Value = 2394
$ws.Range("L129").Value = 3333
$ws.Range("M129").Value = 2606
$ws.Range("N129").Value = -13333

$ws.Range("H131").Value = 2551192.2
$ws.Range("I131").Value = 8257.857
$ws.Range("J131").Value = 3737895
$ws.Range("K131").Value = 24773.571
$ws.Range("L131").Value = 11213685
$ws.Range("M131").Value = -19733.571
$ws.Range("N131").Value = -11223765

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9833.333000000001
$ws.Range("I80").Value = 6000
$ws.Range("J80").Value = 11750
$ws.Range("K80").Value = 6000
$ws.Range("L80").Value = 11750
$ws.Range("M80").Value = -5002
$ws.Range("N80").Value = -13746

$ws.Range("H83").Value = 9833.333000000001
$ws.Range("I83").Value = 6000
$ws.Range("J83").Value = 11750
$ws.Range("K83").Value = 30000
$ws.Range("L83").Value = 58750
$ws.Range("M83").Value = -25008
$ws.Range("N83").Value = -68734

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1199.5
$ws.Range("I82").Value = 1199.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1199.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -838.5
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 1199.5
$ws.Range("I85").Value = 1199.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1199.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 48.5
$ws.Range("N85").ClearContents()

$ws.Range("H132").Value = 23297.541
$ws.Range("I132").Value = 36685.285
$ws.Range("J132").Value = 4554.7
$ws.Range("K132").Value = 110055.855
$ws.Range("L132").Value = 13664.1
$ws.Range("M132").Value = -107525.855
$ws.Range("N132").Value = -18724.1

$ws.Range("H133").Value = 29999
$ws.Range("J133").Value = 29999
$ws.Range("L133").Value = 29999
$ws.Range("N133").Value = -35059

$ws.Range("H136").Value = 11951.5
$ws.Range("I136").Value = 33833.332
$ws.Range("J136").Value = 2573.5715
$ws.Range("K136").Value = 101499.996
$ws.Range("L136").Value = 7720.7145
$ws.Range("M136").Value = -98949.99600000001
$ws.Range("N136").Value = -12820.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1114.6875
$ws.Range("I122").Value = 1083.3334
$ws.Range("K122").Value = 3250.0002
$ws.Range("M122").Value = -800.0001999999999

$ws.Range("H132").Value = 1580.3334
$ws.Range("I132").Value = 1197.8
$ws.Range("J132").Value = 3493
$ws.Range("K132").Value = 3593.4
$ws.Range("L132").Value = 10479
$ws.Range("M132").Value = -1063.4
$ws.Range("N132").Value = -15539

$ws.Range("H136").Value = 6593.45
$ws.Range("I136").Value = 8054.3125
$ws.Range("K136").Value = 24162.9375
$ws.Range("M136").Value = -21612.9375
